$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.164.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.09%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.514.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'321.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'109.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.74%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.30%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'40.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +11.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0819"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.65%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.15%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.913.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.45%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.514.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.80%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.850"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'48.014.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.94%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.91%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.38%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0945"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.83%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.52%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'272.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +11.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.51%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'25.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.23%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.68%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'10.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.34%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.144"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.80%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'35.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'49.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.52%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'19.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.48%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0783"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.36%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.57%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.111"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.62%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'121.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.89%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'21.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.97%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.027.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +7.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'5.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.48%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'79.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.05%  "
$ws.Range("E51").Style = "Normal"
